$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 90 ---
# Column A must be stored as TEXT "89" (not the number 89). Typing the digit
# string directly gets auto-coerced to a number, so build it as a formula
# and convert the formula result to a literal value via copy/paste-special.
$ws.Range("A90").Formula = '="89"'
$ws.Range("A90").Copy()
$ws.Range("A90").PasteSpecial(-4163)

$ws.Range("B90").Value = 43984.0
$ws.Range("B90").NumberFormat = "m/d/yyyy"
$ws.Range("C90").Value = 2286.0
$ws.Range("D90").Value = 183083.0
$ws.Range("E90").Value = 3.0

# --- Row 91 ---
$ws.Range("A91").Formula = '="90"'
$ws.Range("A91").Copy()
$ws.Range("A91").PasteSpecial(-4163)

$ws.Range("B91").Value = 43983.0
$ws.Range("B91").NumberFormat = "m/d/yyyy"
$ws.Range("C91").Value = 2297.0
$ws.Range("D91").Value = 185380.0
$ws.Range("E91").Value = 2.0

$excel.CutCopyMode = $false
